$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1, Col 1: 85÷9= -> 60÷3=
$cell = $t.Cell(1, 1)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "85÷9=") { throw "Unexpected cell text at (1,1): $($cell.Range.Text)" }
$cell.Range.Text = "60÷3="

# Row 1, Col 2: 83÷9= -> 20÷4=
$cell = $t.Cell(1, 2)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "83÷9=") { throw "Unexpected cell text at (1,2): $($cell.Range.Text)" }
$cell.Range.Text = "20÷4="

# Row 1, Col 3: 17÷2= -> 50÷8=
$cell = $t.Cell(1, 3)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "17÷2=") { throw "Unexpected cell text at (1,3): $($cell.Range.Text)" }
$cell.Range.Text = "50÷8="

# Row 1, Col 4: 28÷4= -> 59÷5=
$cell = $t.Cell(1, 4)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "28÷4=") { throw "Unexpected cell text at (1,4): $($cell.Range.Text)" }
$cell.Range.Text = "59÷5="

# Row 1, Col 5: 14÷9= -> 93÷9=
$cell = $t.Cell(1, 5)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "14÷9=") { throw "Unexpected cell text at (1,5): $($cell.Range.Text)" }
$cell.Range.Text = "93÷9="

# Row 5, Col 1: 58÷2= -> 65÷2=
$cell = $t.Cell(5, 1)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "58÷2=") { throw "Unexpected cell text at (5,1): $($cell.Range.Text)" }
$cell.Range.Text = "65÷2="

# Row 5, Col 2: 10÷5= -> 14÷6=
$cell = $t.Cell(5, 2)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "10÷5=") { throw "Unexpected cell text at (5,2): $($cell.Range.Text)" }
$cell.Range.Text = "14÷6="

# Row 5, Col 3: 41÷7= -> 85÷4=
$cell = $t.Cell(5, 3)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "41÷7=") { throw "Unexpected cell text at (5,3): $($cell.Range.Text)" }
$cell.Range.Text = "85÷4="

# Row 5, Col 4: 41÷7= -> 40÷3=
$cell = $t.Cell(5, 4)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "41÷7=") { throw "Unexpected cell text at (5,4): $($cell.Range.Text)" }
$cell.Range.Text = "40÷3="

# Row 5, Col 5: 28÷4= -> 25÷2=
$cell = $t.Cell(5, 5)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "28÷4=") { throw "Unexpected cell text at (5,5): $($cell.Range.Text)" }
$cell.Range.Text = "25÷2="

# Row 9, Col 1: 30÷2= -> 85÷9=
$cell = $t.Cell(9, 1)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "30÷2=") { throw "Unexpected cell text at (9,1): $($cell.Range.Text)" }
$cell.Range.Text = "85÷9="

# Row 9, Col 2: 72÷3= -> 46÷5=
$cell = $t.Cell(9, 2)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "72÷3=") { throw "Unexpected cell text at (9,2): $($cell.Range.Text)" }
$cell.Range.Text = "46÷5="

# Row 9, Col 3: 94÷7= -> 32÷8=
$cell = $t.Cell(9, 3)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "94÷7=") { throw "Unexpected cell text at (9,3): $($cell.Range.Text)" }
$cell.Range.Text = "32÷8="

# Row 9, Col 4: 39÷6= -> 66÷3=
$cell = $t.Cell(9, 4)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "39÷6=") { throw "Unexpected cell text at (9,4): $($cell.Range.Text)" }
$cell.Range.Text = "66÷3="

# Row 9, Col 5: 50÷3= -> 32÷6=
$cell = $t.Cell(9, 5)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "50÷3=") { throw "Unexpected cell text at (9,5): $($cell.Range.Text)" }
$cell.Range.Text = "32÷6="

# Row 13, Col 1: 66÷2= -> 69÷5=
$cell = $t.Cell(13, 1)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "66÷2=") { throw "Unexpected cell text at (13,1): $($cell.Range.Text)" }
$cell.Range.Text = "69÷5="

# Row 13, Col 2: 12÷3= -> 71÷5=
$cell = $t.Cell(13, 2)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "12÷3=") { throw "Unexpected cell text at (13,2): $($cell.Range.Text)" }
$cell.Range.Text = "71÷5="

# Row 13, Col 3: 73÷3= -> 97÷8=
$cell = $t.Cell(13, 3)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "73÷3=") { throw "Unexpected cell text at (13,3): $($cell.Range.Text)" }
$cell.Range.Text = "97÷8="

# Row 13, Col 4: 78÷8= -> 93÷5=
$cell = $t.Cell(13, 4)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "78÷8=") { throw "Unexpected cell text at (13,4): $($cell.Range.Text)" }
$cell.Range.Text = "93÷5="

# Row 13, Col 5: 54÷8= -> 86÷7=
$cell = $t.Cell(13, 5)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "54÷8=") { throw "Unexpected cell text at (13,5): $($cell.Range.Text)" }
$cell.Range.Text = "86÷7="

# Row 17, Col 1: 36÷5= -> 91÷9=
$cell = $t.Cell(17, 1)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "36÷5=") { throw "Unexpected cell text at (17,1): $($cell.Range.Text)" }
$cell.Range.Text = "91÷9="

# Row 17, Col 2: 78÷9= -> 72÷4=
$cell = $t.Cell(17, 2)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "78÷9=") { throw "Unexpected cell text at (17,2): $($cell.Range.Text)" }
$cell.Range.Text = "72÷4="

# Row 17, Col 3: 96÷7= -> 90÷4=
$cell = $t.Cell(17, 3)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "96÷7=") { throw "Unexpected cell text at (17,3): $($cell.Range.Text)" }
$cell.Range.Text = "90÷4="

# Row 17, Col 4: 84÷4= -> 81÷6=
$cell = $t.Cell(17, 4)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "84÷4=") { throw "Unexpected cell text at (17,4): $($cell.Range.Text)" }
$cell.Range.Text = "81÷6="

# Row 17, Col 5: 40÷5= -> 91÷6=
$cell = $t.Cell(17, 5)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "40÷5=") { throw "Unexpected cell text at (17,5): $($cell.Range.Text)" }
$cell.Range.Text = "91÷6="
